# "clarify status report for unit tests"
#
# The PASS/FAIL formula in column D (UnitTests/Tests/TimeSeries.xlsx,
# Sheet1) is extended so a test whose *expected* value (column B) is itself
# an error reports "ERROR" instead of silently comparing two errors and
# calling it a PASS/FAIL. D4:D11 is one shared formula group (si="0"), so
# it is rewritten as a single range assignment (mirrors selecting D4:D11
# and entering the formula with Ctrl+Enter in real Excel) to keep the
# <f t="shared" .../> grouping intact instead of exploding every cell into
# its own literal formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D3 is a standalone (non-shared) formula.
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'

# D4:D11 share formula si="0" (master cell D4). Assigning the whole range
# at once keeps it a single shared-formula group; Excel/this engine
# relativizes B4/C4 per destination row the same way AutoFill would.
$ws.Range("D4:D11").Formula = '=IF(ISERROR(B4),"ERROR",IF(ISERROR(C4),"FAIL",IF(B4=C4,"PASS","FAIL")))'
